# Business partner bind include,exclude validation
# Replace the Productgroup.xlsx sample data: trim the product-group list
# down to two rows (Galaxy / grptest) and fix the header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old data range first so no stray rows/shared strings remain.
$ws.Range("A1:B17").Clear()

# Header row
$ws.Range("A1").Value = "Product Group Name"
$ws.Range("B1").Value = "Level"

# Data rows (write "grptest" first so it lands before "Galaxy" in the
# workbook's shared-string table, matching the authored file).
$ws.Range("A3").Value = "grptest"
$ws.Range("B3").Value = 3

$ws.Range("A2").Value = "Galaxy"
$ws.Range("B2").Value = 232

# Match the saved selection state from the authored workbook.
$ws.Range("C2").Select()
